# "Update countries & provincias Spain"
#
# Refresh of the COVID-19 "Pais" dashboard: the data snapshot moved from
# 11-Apr-2020 15:22 to 15:52, so the timestamp banner and a handful of
# per-country rows (Casos totales / Nuevos casos / Casos activos /
# Recuperados / Casos criticos / Muertes hoy / Muertes) get new numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# --- Timestamp banner (A1) -------------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 11 de Abril de 2020 a las 15:52"

# --- Reino Unido (row 10) ---------------------------------------------------
$ws.Cells.Item(10, 2).Value = 78991
$ws.Cells.Item(10, 3).Value = 5233
$ws.Cells.Item(10, 5).Value = 68772
$ws.Cells.Item(10, 7).Value = 917
$ws.Cells.Item(10, 8).Value = 9875

# --- Islandia (row 59) -------------------------------------------------------
$ws.Cells.Item(59, 2).Value = 1689
$ws.Cells.Item(59, 3).Value = 14
$ws.Cells.Item(59, 4).Value = 841
$ws.Cells.Item(59, 5).Value = 841

# --- Bosnia y Herzegovina (row 74) ------------------------------------------
$ws.Cells.Item(74, 4).Value = 139
$ws.Cells.Item(74, 5).Value = 759

# --- Ghana (row 97) ----------------------------------------------------------
$ws.Cells.Item(97, 5).Value = 396
$ws.Cells.Item(97, 7).Value = 2
$ws.Cells.Item(97, 8).Value = 8

# --- Ruanda (row 129) ---------------------------------------------------------
$ws.Cells.Item(129, 4).Value = 18
$ws.Cells.Item(129, 5).Value = 100

# --- Antigua y Barbuda (row 164) --------------------------------------------
$ws.Cells.Item(164, 4).Value = 1
$ws.Cells.Item(164, 6).Value = 0
$ws.Cells.Item(164, 8).Value = 1

# --- Somalia (row 165) -------------------------------------------------------
$ws.Cells.Item(165, 4).Value = 0
$ws.Cells.Item(165, 6).Value = 1
$ws.Cells.Item(165, 8).Value = 2

# --- Angola (row 167) ---------------------------------------------------------
$ws.Cells.Item(167, 3).Value = 2

# --- Siria (row 168) -----------------------------------------------------------
$ws.Cells.Item(168, 4).Value = 2
$ws.Cells.Item(168, 5).Value = 15

# --- Maldivas (row 169) --------------------------------------------------------
$ws.Cells.Item(169, 4).Value = 4
$ws.Cells.Item(169, 5).Value = 13
$ws.Cells.Item(169, 8).Value = 2

# --- Laos (row 170) -------------------------------------------------------------
$ws.Cells.Item(170, 2).Value = 19
$ws.Cells.Item(170, 3).Value = 0
$ws.Cells.Item(170, 4).Value = 13
$ws.Cells.Item(170, 5).Value = 6

# --- Nueva Caledonia (row 171) --------------------------------------------------
$ws.Cells.Item(171, 3).Value = 2
$ws.Cells.Item(171, 4).Value = 0
$ws.Cells.Item(171, 5).Value = 18

# --- Guinea Ecuatorial (row 172) ------------------------------------------------
$ws.Cells.Item(172, 4).Value = 1
$ws.Cells.Item(172, 5).Value = 17

# --- Islas Virgenes de los Estados Unidos (row 173) -----------------------------
$ws.Cells.Item(173, 2).Value = 18
$ws.Cells.Item(173, 4).Value = 3
$ws.Cells.Item(173, 5).Value = 15

# --- Sudan (row 174) -------------------------------------------------------------
$ws.Cells.Item(174, 4).Value = 0
$ws.Cells.Item(174, 5).Value = 17
$ws.Cells.Item(174, 8).Value = 0

# --- Suazilandia (row 186) --------------------------------------------------------
$ws.Cells.Item(186, 3).Value = 3
$ws.Cells.Item(186, 4).Value = 0
$ws.Cells.Item(186, 5).Value = 10
$ws.Cells.Item(186, 6).Value = 1
$ws.Cells.Item(186, 7).Value = 1
$ws.Cells.Item(186, 8).Value = 2

# --- Seychelles (row 187) -----------------------------------------------------------
$ws.Cells.Item(187, 2).Value = 12
$ws.Cells.Item(187, 4).Value = 7
$ws.Cells.Item(187, 5).Value = 5

# --- Republica del Chad (row 188) ---------------------------------------------------
$ws.Cells.Item(188, 4).Value = 0
$ws.Cells.Item(188, 5).Value = 11

# --- Groenlandia (row 189) -----------------------------------------------------------
$ws.Cells.Item(189, 4).Value = 2
$ws.Cells.Item(189, 5).Value = 9

# --- Belice (row 190) ------------------------------------------------------------------
$ws.Cells.Item(190, 2).Value = 11
$ws.Cells.Item(190, 4).Value = 11
$ws.Cells.Item(190, 5).Value = 0
$ws.Cells.Item(190, 6).Value = 0
$ws.Cells.Item(190, 8).Value = 0

# --- Surinam (row 191) ------------------------------------------------------------------
$ws.Cells.Item(191, 4).Value = 0
$ws.Cells.Item(191, 5).Value = 8
$ws.Cells.Item(191, 6).Value = 1
$ws.Cells.Item(191, 8).Value = 2

# --- Malaui (row 192) --------------------------------------------------------------------
$ws.Cells.Item(192, 2).Value = 10
$ws.Cells.Item(192, 4).Value = 4
$ws.Cells.Item(192, 5).Value = 5
$ws.Cells.Item(192, 6).Value = 0
